# Insert a new data row before row 277 (shifts existing rows 277-343 down to 278-344)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(277).Insert()

# Populate the newly inserted row 277 with the new record's data
$ws.Range("A277").Value = 10
$ws.Range("B277").Value = "Vega Modelo de Temuco"
$ws.Range("C277").Value = "La Araucanía"
$ws.Range("D277").Value = 45015
$ws.Range("E277").Value = 9
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100103
$ws.Range("H277").Value = "Frutos de hueso (carozo)"
$ws.Range("I277").Value = 100103002
$ws.Range("J277").Value = "Ciruela"
$ws.Range("K277").Value = "Pink Delight"
$ws.Range("L277").Value = "Primera"
$ws.Range("M277").Value = 85
$ws.Range("N277").Value = 14000
$ws.Range("O277").Value = 15000
$ws.Range("P277").Value = 14647
$ws.Range("Q277").Value = "$/bandeja 18 kilos granel"
$ws.Range("R277").Value = "Región de O'Higgins"
$ws.Range("S277").Value = 814
$ws.Range("T277").Value = 18
